$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.165.70'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '2.306.58'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.60'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.510'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.75%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.69'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.55%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.71'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.53%  '
$ws.Range('E14').Value = '  +2.14%  '
$ws.Range('D15').Value = '2.665.39'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').Value = '2.319.45'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('D18').Value = '43.078.61'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('E19').Value = '  +9.53%  '
$ws.Range('D20').Value = '0.0₃0906'
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.98'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.32'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.28%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('E34').Value = '  +4.69%  '
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0690'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('E38').Value = '  +1.50%  '
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.81'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('E42').Value = '  +3.14%  '
$ws.Range('D43').Value = '1.979.65'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.28'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.84'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.35'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.70%  '
$ws.Range('E49').Value = '  +4.04%  '
$ws.Range('D50').Value = '2.531.41'
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.91'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.24%  '
